$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the stray top row ("A","B","C") and shift everything up, which
# also drops the now-unused last row (row 12).
$ws.Rows.Item(1).Delete()

# Reset the saved selection back to the default top-left cell (matches a
# freshly-reopened sheet instead of the old, now-stale C7 selection).
[void]$ws.Range("A1").Select()
